# Fixed a bug in Mask
# Re-orders the data rows (A2:F25) of Sheet1 to the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newData = @(
    @(902,1,0,0,0,0),
    @(1201,2,10,10,10,10),
    @(401,9,48,67,75,45),
    @(601,9,60,67,60,42),
    @(501,9,52,30,75,45),
    @(101,9,30,15,60,15),
    @(901,16,15,45,60,60),
    @(1001,18,30,75,60,72),
    @(301,6,45,30,60,45),
    @(1203,3,15,15,15,15),
    @(701,3,90,45,97,15),
    @(201,9,30,15,45,30),
    @(801,3,67,65,52,45),
    @(1202,2,10,10,10,10),
    @(2,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(1101,0,15,30,30,0),
    @(502,0,4,0,0,0),
    @(802,0,4,5,4,0),
    @(1,0,2,2,2,2),
    @(402,0,0,4,0,0),
    @(602,0,0,4,0,9),
    @(702,0,0,0,4,0),
    @(1002,0,0,0,0,9)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newData[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $rowValues[$col]
    }
}
